# Remove the "Palestine / Masjid-Al-Aqsa" row (row 66) from the data table.
# Deleting the entire row shifts all following rows up by one, which matches
# the target diff (old row 67 becomes new row 66, ..., old row 105 is gone,
# and the sheet dimension shrinks from A1:E105 to A1:E104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(66).Delete()
